# Weekly update: insert a new daily record at the top of the
# "Vega Modelo de Temuco - Melón" data block (row 732), shifting all
# subsequent rows down by one.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at row 732; existing rows 732:849 shift to 733:850
$ws.Rows.Item(732).Insert()

# Populate the new row 732 with the new day's record
$ws.Cells.Item(732, 1).Value = 10
$ws.Cells.Item(732, 2).Value = "Vega Modelo de Temuco"
$ws.Cells.Item(732, 3).Value = "La Araucanía"
$ws.Cells.Item(732, 4).Value = 45218
$ws.Cells.Item(732, 5).Value = 9
$ws.Cells.Item(732, 6).Value = 100112027
$ws.Cells.Item(732, 7).Value = "Melón"
$ws.Cells.Item(732, 8).Value = "Tuna"
$ws.Cells.Item(732, 9).Value = "Primera"
$ws.Cells.Item(732, 10).Value = 320
$ws.Cells.Item(732, 11).Value = 1800
$ws.Cells.Item(732, 12).Value = 1800
$ws.Cells.Item(732, 13).Value = 1800
$ws.Cells.Item(732, 14).Value = "$/unidad"
$ws.Cells.Item(732, 15).Value = "Perú"
$ws.Cells.Item(732, 16).Value = 1800
$ws.Cells.Item(732, 17).Value = 1
$ws.Cells.Item(732, 18).Value = "Hortaliza"
